$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 865.3333
$ws.Range("I19").Value = 1499
$ws.Range("J19").Value = 548.5
$ws.Range("K19").Value = 1499
$ws.Range("L19").Value = 548.5
$ws.Range("M19").Value = -1324
$ws.Range("N19").Value = -898.5
$ws.Range("H93").Value = 49999.5
$ws.Range("J93").Value = 49999.5
$ws.Range("L93").Value = 49999.5
$ws.Range("N93").Value = -54991.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4698.4614
$ws.Range("I32").Value = 4990
$ws.Range("K32").Value = 4990
$ws.Range("M32").Value = -4703
$ws.Range("H61").Value = 3250
$ws.Range("I61").Value = 3000
$ws.Range("J61").Value = 3500
$ws.Range("K61").Value = 3000
$ws.Range("L61").Value = 3500
$ws.Range("M61").Value = -2788
$ws.Range("N61").Value = -3924
$ws.Range("H80").Value = 30776.75
$ws.Range("J80").Value = 85110
$ws.Range("L80").Value = 85110
$ws.Range("N80").Value = -87106
$ws.Range("H83").Value = 30776.75
$ws.Range("J83").Value = 85110
$ws.Range("L83").Value = 255330
$ws.Range("N83").Value = -265314
$ws.Range("H96").Value = 67500
$ws.Range("J96").Value = 67500
$ws.Range("L96").Value = 67500
$ws.Range("N96").Value = -72992
$ws.Range("H102").Value = 3121
$ws.Range("I102").Value = 3121
$ws.Range("K102").Value = 3121
$ws.Range("M102").Value = -1499
$ws.Range("H110").Value = 1445.5714
$ws.Range("I110").Value = 1445.5714
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1445.5714
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 599.4286
$ws.Range("N110").Value = ""
$ws.Range("H136").Value = 3250
$ws.Range("I136").Value = 3000
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 9000
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -6450
$ws.Range("N136").Value = -15600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4003
$ws.Range("I134").Value = 3670.6667
$ws.Range("K134").Value = 11012.0001
$ws.Range("M134").Value = -8477.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = ""
$ws.Range("N38").Value = ""
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = ""
$ws.Range("N46").Value = ""
$ws.Range("H55").Value = 7536.5
$ws.Range("I55").Value = 10073
$ws.Range("J55").Value = 5000
$ws.Range("K55").Value = 10073
$ws.Range("L55").Value = 5000
$ws.Range("M55").Value = -9758
$ws.Range("N55").Value = -5630
$ws.Range("H58").Value = 3500
$ws.Range("I58").Value = 2000
$ws.Range("K58").Value = 2000
$ws.Range("M58").Value = -1797
$ws.Range("H86").Value = 55000
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = ""
$ws.Range("H89").Value = 55000
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = ""
$ws.Range("H136").Value = 3500
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -3450

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 9766.666999999999
$ws.Range("J55").Value = 13450
$ws.Range("L55").Value = 40350
$ws.Range("N55").Value = -40704
$ws.Range("H81").Value = 4020.4167
$ws.Range("I81").Value = 2849.8
$ws.Range("J81").Value = 4856.5713
$ws.Range("K81").Value = 8549.400000000001
$ws.Range("L81").Value = 14569.7139
$ws.Range("M81").Value = -7426.400000000001
$ws.Range("N81").Value = -16815.7139
$ws.Range("H84").Value = 4020.4167
$ws.Range("I84").Value = 2849.8
$ws.Range("J84").Value = 4856.5713
$ws.Range("K84").Value = 25648.2
$ws.Range("L84").Value = 43709.14169999999
$ws.Range("M84").Value = -20032.2
$ws.Range("N84").Value = -54941.14169999999
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 93.166664
$ws.Range("I2").Value = 89.75
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 89.75
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = 23.25
$ws.Range("N2").Value = -326
$ws.Range("H11").Value = 625123.75
$ws.Range("I11").Value = 1000000
$ws.Range("J11").Value = 250247.5
$ws.Range("K11").Value = 1000000
$ws.Range("L11").Value = 250247.5
$ws.Range("M11").Value = -999861
$ws.Range("N11").Value = -250525.5
$ws.Range("H43").Value = 30526
$ws.Range("I43").Value = 31578
$ws.Range("J43").Value = 30000
$ws.Range("K43").Value = 31578
$ws.Range("L43").Value = 30000
$ws.Range("M43").Value = -31427
$ws.Range("N43").Value = -30302
$ws.Range("H57").Value = 10675
$ws.Range("I57").Value = 10675
$ws.Range("K57").Value = 10675
$ws.Range("M57").Value = -9855
$ws.Range("H80").Value = 22500
$ws.Range("I80").Value = 15000
$ws.Range("K80").Value = 15000
$ws.Range("M80").Value = -14002
$ws.Range("H83").Value = 22500
$ws.Range("I83").Value = 15000
$ws.Range("K83").Value = 75000
$ws.Range("M83").Value = -70008

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1000.1
$ws.Range("I16").Value = 1086.8334
$ws.Range("K16").Value = 1086.8334
$ws.Range("M16").Value = -916.8334
$ws.Range("H132").Value = 4773.8335
$ws.Range("J132").Value = 5666
$ws.Range("L132").Value = 16998
$ws.Range("N132").Value = -22058

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").Value = ""
